$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell values for the existing rows (A1:B4) and add a new row (A5:B5).
$ws.Range("A1").Value = -0.010067499149145369
$ws.Range("B1").Value = 0.010067498034236068
$ws.Range("A2").Value = 0.044028105706899627
$ws.Range("B2").Value = -0.044028106842829422
$ws.Range("A3").Value = -0.068461654260332142
$ws.Range("B3").Value = 0.068461653146266982
$ws.Range("A4").Value = 0.00088045523342019695
$ws.Range("B4").Value = -0.00088045649764036614
$ws.Range("A5").Value = -0.0027517693972807359
$ws.Range("B5").Value = 0.0027517680889708749

# Widen columns A and B (closest attainable widths on this host's column-width grid).
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666
